$d = $word.ActiveDocument

$pairs = @(
    @{old = "437×2="; new = "887×8="},
    @{old = "639×3="; new = "953×5="},
    @{old = "729×8="; new = "463×4="},
    @{old = "626×8="; new = "141×4="},
    @{old = "222×2="; new = "579×5="},
    @{old = "548×2="; new = "169×5="},
    @{old = "597×6="; new = "965×7="},
    @{old = "788×5="; new = "407×4="},
    @{old = "342×9="; new = "195×8="},
    @{old = "772×3="; new = "186×4="},
    @{old = "734×6="; new = "807×7="},
    @{old = "225×3="; new = "479×9="},
    @{old = "741×4="; new = "566×8="},
    @{old = "621×4="; new = "857×4="},
    @{old = "383×5="; new = "716×8="},
    @{old = "967×6="; new = "794×5="},
    @{old = "306×9="; new = "906×8="},
    @{old = "554×2="; new = "887×5="},
    @{old = "261×6="; new = "616×2="},
    @{old = "984×7="; new = "868×9="},
    @{old = "904×7="; new = "944×5="},
    @{old = "445×6="; new = "632×9="},
    @{old = "511×3="; new = "133×7="},
    @{old = "594×9="; new = "859×2="},
    @{old = "473×2="; new = "407×5="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $pair.new, 2)
}

$d.Save()
